$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Informe-01-010073-A-TC-TM-TP")

# Column D = municipio-nombre: it was wrongly documented as a measure; now
# it is documented as the proper sdmx dimension for the area, with its
# dim/measure flag and URI type updated accordingly.
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# Column F = tipo-licencias-descripcion: the inverse move, it now becomes
# a measure instead of a dimension.
$ws.Range("F2").Value = "iaest-measure:tipo-licencias-descripcion"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"

# The mapping file reference for tipo-licencias-descripcion no longer
# applies now that it is a measure, not a dimension.
$ws.Range("F5").Clear()
